# Applies the changes described by the commit "break out stock.yaml completed":
#  1. On sheet "day", rows 908-911 column D held the BSE scrip code as text
#     (inlineStr). Convert those four cells to genuine numbers.
#  2. On sheet "week", 17 new rows (558-574) of scraped data were appended,
#     extending the sheet's dimension from A1:I557 to A1:I574.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Part 1: sheet "day" -- convert D908:D911 from text to numeric values
# ---------------------------------------------------------------------
$dayWs = $wb.Worksheets.Item("day")

$dayWs.Cells.Item(908, 4).Value = 500387
$dayWs.Cells.Item(909, 4).Value = 500575
$dayWs.Cells.Item(910, 4).Value = 532755
$dayWs.Cells.Item(911, 4).Value = 532523

# ---------------------------------------------------------------------
# Part 2: sheet "week" -- append new rows 558-574
# ---------------------------------------------------------------------
$weekWs = $wb.Worksheets.Item("week")

# Columns: row, A rank, B symbol, C company name, D scrip code (text),
#          E pct change, F value, G quantity, H period label, I timestamp
$newRows = @(
    @(558, 1, "BOSCHLTD", "Bosch Limited", "500530", 1.96, 33862.1, 25998, "week", "15/11/2024 11:33:13"),
    @(559, 2, "NAUKRI", "Info Edge (india) Limited", "532777", 2.01, 7768.2, 187611, "week", "15/11/2024 11:33:13"),
    @(560, 3, "ATUL", "Atul Limited", "500027", 1.08, 7303.2, 46953, "week", "15/11/2024 11:33:13"),
    @(561, 4, "HDFCAMC", "HDFC Asset Management Company Ltd", "541729", 0.7, 4278.75, 393645, "week", "15/11/2024 11:33:13"),
    @(562, 5, "TORNTPHARM", "Torrent Pharmaceuticals Limited", "500420", -0.46, 3102.4, 223405, "week", "15/11/2024 11:33:13"),
    @(563, 6, "LALPATHLAB", "Dr. Lal Path Labs Ltd.", "539524", -0.03, 2956.4, 149233, "week", "15/11/2024 11:33:13"),
    @(564, 7, "ADANIENT", "Adani Enterprises Limited", "512599", 0.36, 2826.8, 646797, "week", "15/11/2024 11:33:13"),
    @(565, 8, "METROPOLIS", "Metropolis Healthcare Ltd", "542650", -0.88, 2032.2, 99058, "week", "15/11/2024 11:33:13"),
    @(566, 9, "OBEROIRLTY", "Oberoi Realty Limited", "533273", 2.76, 1980.3, 1392570, "week", "15/11/2024 11:33:13"),
    @(567, 10, "COROMANDEL", "Coromandel International Limited", "506395", 1.84, 1715.95, 572592, "week", "15/11/2024 11:33:13"),
    @(568, 11, "HDFCBANK", "Hdfc Bank Limited", "500180", 0.68, 1692.75, 13573235, "week", "15/11/2024 11:33:13"),
    @(569, 12, "CIPLA", "Cipla Limited", "500087", -0.39, 1499.75, 1469143, "week", "15/11/2024 11:33:13"),
    @(570, 13, "PEL", "Piramal Enterprises Limited", "500302", 3.16, 1044.25, 961682, "week", "15/11/2024 11:33:13"),
    @(571, 14, "SBIN", "State Bank Of India", "500112", -0.54, 804.25, 9461484, "week", "15/11/2024 11:33:13"),
    @(572, 15, "INDHOTEL", "The Indian Hotels Company Limited", "500850", 3.81, 741.35, 5837816, "week", "15/11/2024 11:33:13"),
    @(573, 16, "RECLTD", "Rural Electrification Corporation Limited", "532955", -1.2, 502.35, 7471948, "week", "15/11/2024 11:33:13"),
    @(574, 17, "LAURUSLABS", "Laurus Labs Limited", "540222", 1.96, 486, 878316, "week", "15/11/2024 11:33:13")
)

# Force column D to text formatting for the new block so the scrip codes
# (e.g. "500530") are stored as strings, not auto-converted to numbers.
$dCodeRange = $weekWs.Range("D558:D574")
$dCodeRange.NumberFormat = "@"

foreach ($row in $newRows) {
    $r = $row[0]
    $weekWs.Cells.Item($r, 1).Value = $row[1]
    $weekWs.Cells.Item($r, 2).Value = $row[2]
    $weekWs.Cells.Item($r, 3).Value = $row[3]
    $weekWs.Cells.Item($r, 4).Value = $row[4]
    $weekWs.Cells.Item($r, 5).Value = $row[5]
    $weekWs.Cells.Item($r, 6).Value = $row[6]
    $weekWs.Cells.Item($r, 7).Value = $row[7]
    $weekWs.Cells.Item($r, 8).Value = $row[8]
    $weekWs.Cells.Item($r, 9).Value = $row[9]
}

# Clear the temporary style so the new cells end up with the sheet's
# default (unstyled) formatting, matching the rest of the data.
$dCodeRange.Style = "Normal"

Write-Output "edit complete"
